$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.001.58"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.642.31"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "215.88"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.870.71"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.667.74"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "0.0₃0762"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "63.41"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "26.093.01"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "194.80"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "6.21"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "1.80"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "143.19"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "1.129.41"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "99.24"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "0.795"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "1.779.50"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").Value = "56.64"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").Value = "7.73"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.28%  "
